$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append year data rows 199-210 (A: 197-208, B: oil price values)
# mirrors the existing pattern in column A (s="1" style) and plain column B

$data = @(
    @(197, 0.7800829875518671),
    @(198, 0.4835910976989815),
    @(199, 0.546058091286307),
    @(200, 0.5062240663900415),
    @(201, 0.313278008298755),
    @(202, 0.6484884410195613),
    @(203, 0.5767634854771784),
    @(204, 0.6979253112033195),
    @(205, 0.6556016597510372),
    @(206, 0.6224066390041493),
    @(207, 0.7053941908713692),
    @(208, 0.5560165975103734)
)

$startRow = 199
$templateRow = 198
$templateCell = $ws.Cells.Item($templateRow, 1)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $aVal = $data[$i][0]
    $bVal = $data[$i][1]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $aVal

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $bVal
}

# Copy the formatting (style index) of the last existing "A" template cell
# down across the newly appended A column cells, matching the original
# pattern used throughout the sheet.
$templateCell.Copy() | Out-Null
$destRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $data.Length - 1, 1))
$destRange.PasteSpecial(-4122) | Out-Null
